$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation so numeric-looking strings (e.g. "1.00", "26.659.37")
# are preserved exactly as literal text instead of being parsed into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.708.52"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "1.645.77"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "215.86"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "0.507"
$ws.Range("E6").Value = "  +1.62%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("D9").Value = "0.0626"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "19.16"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.874.96"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "1.633.68"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "4.18"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "65.18"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "26.701.01"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "0.0₃0743"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "217.91"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "4.36"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").Value = "6.27"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "9.51"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("D24").Value = "2.26"
$ws.Range("E24").Value = "  +12.74%  "
$ws.Range("D25").Value = "145.86"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "7.12"
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("D29").Value = "15.72"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").Value = "0.0516"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").Value = "3.01"
$ws.Range("E33").Value = "  +2.39%  "
$ws.Range("D34").Value = "1.278.67"
$ws.Range("E34").Value = "  +4.83%  "
$ws.Range("E35").Value = "  +3.59%  "
$ws.Range("E36").Value = "  +1.54%  "
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("D38").Value = "0.535"
$ws.Range("E38").Value = "  +5.93%  "
$ws.Range("D39").Value = "0.825"
$ws.Range("E39").Value = "  +4.04%  "
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("D41").Value = "0.813"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").Value = "5.45"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").Value = "1.784.72"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").Value = "91.93"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("D46").Value = "59.90"
$ws.Range("E46").Value = "  +9.21%  "
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0515"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.75"
$ws.Range("E50").Value = "  +3.32%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0968"
$ws.Range("E51").Value = "  +1.98%  "

# Remove the temporary text-number-format override so the cell style
# reverts to the workbook default (matches original formatting).
$ws.Range("D2:E51").ClearFormats()

